$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 0, 1.104883657715537),
    @(3, 0.6545652718822623, 1.626987699542094,  0.1496068669990043, 0.5333859586016987, 0, 2.964545797025059),
    @(4, 0.6545652718822623, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 0, 1.642449346116345),
    @(5, 0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1, 1.104883657715537),
    @(6, 1.445647641019636,  1.626987699542094,  0.7210945179870265, 0.5333859586016987, 0, 4.327115817150455),
    @(7, 3.272327238179451,  1.626987699542094,  3.223369029078222,  0.5333859586016987, 0, 8.656069925401464),
    @(8, 1.445647641019636,  1.626987699542094,  0.1496068669990043, 13.86384647080068,  0, 17.08608867836142),
    @(9, 3.272327238179451,  1.626987699542094,  3.223369029078222,  0.5333859586016987, 1, 8.656069925401464),
    @(10, 0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 3.536033448013082)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
